$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 395 ("こわれている" / "broken") entirely; all rows below shift up by one.
$ws.Rows(395).Delete()
